$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Replace-TextAfter($anchorText, $old, $new) {
    $full = $d.Content
    $anchor = $d.Content
    $anchor.Find.Execute($anchorText) | Out-Null
    $startPos = $anchor.End
    $scoped = $d.Range($startPos, $full.End)
    $scoped.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1. Merge the "Members:" names into a single plain-text run (drops the
#    spell-check proofErr wrappers around "Ringdahl" / "Barua").
Replace-Text "Nicklas Ringdahl, Samir Barua, Mime Liu" "Nicklas Ringdahl, Samir Barua, Mime Liu"

# 2. Hypothesis bullet text changes (crime/covid project pivot).
Replace-Text "People are less compliant with the lockdown rules in this lockdown compared to last year;" "Certain crimes have decreased due to the lockdown (Breaking and entering for example) while others will have increased (Domestic violence as example) "

Replace-Text "When there are more lockdown breaches, there are higher daily case increase rates;" "Areas with higher crime rates will overlap with areas showing high number of covid cases;"

Replace-Text "The outbreak hotspots (types and location of venues) are generally similar in both lockdowns;" "Crime will have dropped early in lockdown and then slowly increased throughout the year;"

Replace-Text "There might be a strong correlation between lockdown rules adherence and the speed of Covid spread." "There will be a big shift in locations of crimes from outdoors to domestic situations."

# 3. Remove the proofErr/gramStart split in "How will you access the data?"
Replace-Text "How will you access the data? (e.g. via an API or downloading the data)" "How will you access the data? (e.g. via an API or downloading the data)"

# 4. Merge "N - Crime data (types, amount etc)" into one run -- only the
#    Phase 2 occurrence (the Phase 1 occurrence has a different run split
#    and must stay untouched).
Replace-TextAfter "Phase 2 - Data processing" "N - Crime data (types, amount etc)" "N - Crime data (types, amount etc)"

# 5. Merge "M – Daily case changes and crime data" into one run (case
#    sensitive match distinguishes it from the lowercase "daily" instance
#    in Phase 3, which must stay untouched).
$d.Content.Find.Execute("M – Daily case changes and crime data", $true, $false, $false, $false, $false, $true, 1, $false, "M – Daily case changes and crime data", 2) | Out-Null
